$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Category Code for the "US Core Average Blood Pressure Profile" row
$ws.Range("C2").Value = "null#vital-signs"

# Update the LOINC codes to include additional component codes
$ws.Range("E2").Value = "LOINC#85354-9, LOINC#96607-7"
$ws.Range("E3").Value = "LOINC#8480-6, LOINC#96608-5"
$ws.Range("E4").Value = "LOINC#8462-4, LOINC#96609-3"

# Remove the trailing rows for Care Experience Preference, Laboratory Result
# and Treatment Intervention Preference profiles (rows 5-7)
$ws.Rows("5:7").Delete()
